$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments (closest achievable values; the engine
#     quantizes ColumnWidth to whole-pixel character widths, same as Excel) ---
$ws.Columns.Item(1).ColumnWidth = 19.833333333333332
$ws.Columns.Item(15).ColumnWidth = 21.833333333333332

# --- Updated ratio/error values (corrected age-calculation algorithm: switched
#     outlier-detection from median- to mean-based) ---
$ws.Range("A2").Value = -2.057741233326138
$ws.Range("B2").Value = 0.00150413777460899
$ws.Range("O2").Value = 1.320948087526761
$ws.Range("P2").Value = 0.3060073134943528
$ws.Range("S2").Value = 0.235023292489963
$ws.Range("T2").Value = 0.5173214114714446
$ws.Range("A3").Value = 9.844867381475719
$ws.Range("B3").Value = 0.0006652207508307103
$ws.Range("O3").Value = 0.003753217049396492
$ws.Range("P3").Value = 1.625466882474972
$ws.Range("S3").Value = 0.001473786502812541
$ws.Range("T3").Value = 1.759794056786792
$ws.Range("A4").Value = -3.816859799880779
$ws.Range("B4").Value = 0.001530490917323254
$ws.Range("O4").Value = 1.321189907788415
$ws.Range("P4").Value = 0.2966429908788535
$ws.Range("S4").Value = 0.2361713386449745
$ws.Range("T4").Value = 0.3786565527152522
$ws.Range("A5").Value = 8.42427471483842
$ws.Range("B5").Value = 0.0008485227401331792
$ws.Range("O5").Value = 0.00799901972189606
$ws.Range("P5").Value = 0.965296865843396
$ws.Range("S5").Value = 0.001540222474305018
$ws.Range("T5").Value = 1.012146500227368
$ws.Range("A6").Value = -0.4662230626985009
$ws.Range("B6").Value = 0.001490489719275895
$ws.Range("O6").Value = 1.320107113994624
$ws.Range("P6").Value = 0.2854457918826794
$ws.Range("S6").Value = 0.2351813653888696
$ws.Range("T6").Value = 0.4348627137034843
$ws.Range("A7").Value = 9.481869175683277
$ws.Range("B7").Value = 0.0006593942755419525
$ws.Range("O7").Value = 0.008774754252940668
$ws.Range("P7").Value = 1.005239137179149
$ws.Range("S7").Value = 0.001654442344015474
$ws.Range("T7").Value = 1.081828161608998
$ws.Range("A8").Value = -1.243280102740218
$ws.Range("B8").Value = 0.001999729052424219
$ws.Range("O8").Value = 1.321383983749202
$ws.Range("P8").Value = 0.3291428443811557
$ws.Range("S8").Value = 0.2370734835371542
$ws.Range("T8").Value = 0.6737925747004632
$ws.Range("A9").Value = 10.93997355051668
$ws.Range("B9").Value = 0.0006244203294787541
$ws.Range("O9").Value = 0.01115021763499328
$ws.Range("P9").Value = 0.7943269013885956
$ws.Range("S9").Value = 0.003289580199064019
$ws.Range("T9").Value = 0.7624518867062802
$ws.Range("A10").Value = -1.393535923240985
$ws.Range("B10").Value = 0.001708883996364041
$ws.Range("O10").Value = 1.321932283906519
$ws.Range("P10").Value = 0.2869589745090427
$ws.Range("S10").Value = 0.2352072486112528
$ws.Range("T10").Value = 0.5397775158218402
$ws.Range("A11").Value = 13.49983411436106
$ws.Range("B11").Value = 0.001705999406143795
$ws.Range("O11").Value = 0.006694459231116382
$ws.Range("P11").Value = 2.096952738741243
$ws.Range("S11").Value = 0.0001438878040375252
$ws.Range("T11").Value = 1.962352159255405
$ws.Range("A12").Value = -3.476725305623707
$ws.Range("B12").Value = 0.001978094688057611
$ws.Range("O12").Value = 1.321976629685961
$ws.Range("P12").Value = 0.2447238217166258
$ws.Range("S12").Value = 0.2346936720085404
$ws.Range("T12").Value = 0.6340033988648147
$ws.Range("A13").Value = 14.19381462357516
$ws.Range("B13").Value = 0.001177851452010102
$ws.Range("O13").Value = 0.0083362536170602
$ws.Range("P13").Value = 1.183913183762592
$ws.Range("S13").Value = 0.002111648560700679
$ws.Range("T13").Value = 1.045225837973764
$ws.Range("A14").Value = -1.483833060055839
$ws.Range("B14").Value = 0.001356315205210567
$ws.Range("O14").Value = 1.325104018511275
$ws.Range("P14").Value = 0.3026695396122626
$ws.Range("S14").Value = 0.2339811975022381
$ws.Range("T14").Value = 0.6475094355142329
$ws.Range("A15").Value = 7.734617743798911
$ws.Range("B15").Value = 0.0006264330546849789
$ws.Range("O15").Value = 0.009683986746654751
$ws.Range("P15").Value = 0.9556486683684929
$ws.Range("S15").Value = 0.0003282266960048436
$ws.Range("T15").Value = 0.9193255083098189
